# Auto-generated Excel COM-interop script
# Applies the Famfrit_Profits market-data refresh: updates currentAveragePrice (H),
# currentAveragePriceNQ (I) and currentAveragePriceHQ (J) for specific leve rows, and
# recomputes the dependent LevePriceNQ (K=I*F), LevePriceHQ (L=J*F), LeveProfitNQ (M=E-K)
# and LeveProfitHQ (N=-2*E-L) columns -- clearing M/N entirely when the corresponding
# K/L value is zero, matching this workbook's existing convention.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 26538.8
$ws.Range("I76").Value = 31431.25
$ws.Range("J76").Value = 6969
$ws.Range("K76").Value = 31431.25
$ws.Range("L76").Value = 6969
$ws.Range("M76").Value = -31116.25
$ws.Range("N76").Value = -7599

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 26538.8
$ws.Range("I79").Value = 31431.25
$ws.Range("J79").Value = 6969
$ws.Range("K79").Value = 31431.25
$ws.Range("L79").Value = 6969
$ws.Range("M79").Value = -30339.25
$ws.Range("N79").Value = -9153

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 648.75
$ws.Range("I92").Value = 518.04
$ws.Range("J92").Value = 1115.5714
$ws.Range("K92").Value = 518.04
$ws.Range("L92").Value = 1115.5714
$ws.Range("M92").Value = 729.96
$ws.Range("N92").Value = -3611.5714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 635.4211
$ws.Range("I98").Value = 635.4211
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 635.4211
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 862.5789
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 635.4211
$ws.Range("I122").Value = 635.4211
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1906.2633
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 543.7366999999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2323.9656
$ws.Range("I122").Value = 2015.2307
$ws.Range("J122").Value = 4999.6665
$ws.Range("K122").Value = 6045.6921
$ws.Range("L122").Value = 14998.9995
$ws.Range("M122").Value = -3595.6921
$ws.Range("N122").Value = -19898.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 72998.8
$ws.Range("I132").Value = 4843.793
$ws.Range("J132").Value = 402414.66
$ws.Range("K132").Value = 14531.379
$ws.Range("L132").Value = 1207243.98
$ws.Range("M132").Value = -12001.379
$ws.Range("N132").Value = -1212303.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10344.4
$ws.Range("I105").Value = 14695.7
$ws.Range("J105").Value = 7443.533
$ws.Range("K105").Value = 14695.7
$ws.Range("L105").Value = 7443.533
$ws.Range("M105").Value = -12948.7
$ws.Range("N105").Value = -10937.533

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 3958.8
$ws.Range("I113").Value = 3958.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3958.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1788.8
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3040.8914
$ws.Range("I31").Value = 1303.963
$ws.Range("J31").Value = 5509.1577
$ws.Range("K31").Value = 1303.963
$ws.Range("L31").Value = 5509.1577
$ws.Range("M31").Value = -1008.963
$ws.Range("N31").Value = -6099.1577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3040.8914
$ws.Range("I34").Value = 1303.963
$ws.Range("J34").Value = 5509.1577
$ws.Range("K34").Value = 1303.963
$ws.Range("L34").Value = 5509.1577
$ws.Range("M34").Value = -1101.963
$ws.Range("N34").Value = -5913.1577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2853.6316
$ws.Range("I58").Value = 3775.3
$ws.Range("J58").Value = 1829.5555
$ws.Range("K58").Value = 3775.3
$ws.Range("L58").Value = 1829.5555
$ws.Range("M58").Value = -3572.3
$ws.Range("N58").Value = -2235.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3523.1667
$ws.Range("I62").Value = 3470.7693
$ws.Range("J62").Value = 3659.4
$ws.Range("K62").Value = 3470.7693
$ws.Range("L62").Value = 3659.4
$ws.Range("M62").Value = -2846.7693
$ws.Range("N62").Value = -4907.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3523.1667
$ws.Range("I65").Value = 3470.7693
$ws.Range("J65").Value = 3659.4
$ws.Range("K65").Value = 17353.8465
$ws.Range("L65").Value = 18297
$ws.Range("M65").Value = -14233.8465
$ws.Range("N65").Value = -24537

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5906.8184
$ws.Range("I99").Value = 6916.6665
$ws.Range("J99").Value = 4695
$ws.Range("K99").Value = 6916.6665
$ws.Range("L99").Value = 4695
$ws.Range("M99").Value = -5418.6665
$ws.Range("N99").Value = -7691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5906.8184
$ws.Range("I126").Value = 6916.6665
$ws.Range("J126").Value = 4695
$ws.Range("K126").Value = 20749.9995
$ws.Range("L126").Value = 14085
$ws.Range("M126").Value = -18279.9995
$ws.Range("N126").Value = -19025

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2853.6316
$ws.Range("I136").Value = 3775.3
$ws.Range("J136").Value = 1829.5555
$ws.Range("K136").Value = 11325.9
$ws.Range("L136").Value = 5488.666499999999
$ws.Range("M136").Value = -8775.900000000001
$ws.Range("N136").Value = -10588.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45.6
$ws.Range("I2").Value = 51.8
$ws.Range("J2").Value = 33.2
$ws.Range("K2").Value = 310.8
$ws.Range("L2").Value = 199.2
$ws.Range("M2").Value = -197.8
$ws.Range("N2").Value = -425.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3334470.8
$ws.Range("I81").Value = 3334470.8
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10003412.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -10002289.4
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 3334470.8
$ws.Range("I84").Value = 3334470.8
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 30010237.2
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -30004621.2
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 95.5
$ws.Range("I2").Value = 40.083332
$ws.Range("J2").Value = 261.75
$ws.Range("K2").Value = 40.083332
$ws.Range("L2").Value = 261.75
$ws.Range("M2").Value = 72.916668
$ws.Range("N2").Value = -487.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4997.8
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 4997.25
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 14991.75
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -19931.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 940.6129
$ws.Range("I16").Value = 894.9643
$ws.Range("J16").Value = 1366.6666
$ws.Range("K16").Value = 894.9643
$ws.Range("L16").Value = 1366.6666
$ws.Range("M16").Value = -724.9643
$ws.Range("N16").Value = -1706.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8335447
$ws.Range("I40").Value = 9805525
$ws.Range("J40").Value = 5004.3335
$ws.Range("K40").Value = 9805525
$ws.Range("L40").Value = 5004.3335
$ws.Range("M40").Value = -9805389
$ws.Range("N40").Value = -5276.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1025.8636
$ws.Range("I55").Value = 672.5454999999999
$ws.Range("J55").Value = 1379.1818
$ws.Range("K55").Value = 672.5454999999999
$ws.Range("L55").Value = 1379.1818
$ws.Range("M55").Value = -499.5454999999999
$ws.Range("N55").Value = -1725.1818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 22000
$ws.Range("I103").Value = 22000
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 22000
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -20828
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2738.261
$ws.Range("I136").Value = 2590.2273
$ws.Range("J136").Value = 5995
$ws.Range("K136").Value = 7770.6819
$ws.Range("L136").Value = 17985
$ws.Range("M136").Value = -5220.6819
$ws.Range("N136").Value = -23085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3750
$ws.Range("I17").Value = 5000
$ws.Range("J17").Value = 2500
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = -4828
$ws.Range("N17").Value = -2844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 3579.889
$ws.Range("I18").Value = 1895.2222
$ws.Range("J18").Value = 5264.5557
$ws.Range("K18").Value = 1895.2222
$ws.Range("L18").Value = 5264.5557
$ws.Range("M18").Value = -1722.2222
$ws.Range("N18").Value = -5610.5557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 6000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 6000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 6000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -6348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 30304134
$ws.Range("I100").Value = 45455708
$ws.Range("J100").Value = 981.2727
$ws.Range("K100").Value = 90911416
$ws.Range("L100").Value = 1962.5454
$ws.Range("M100").Value = -90910875
$ws.Range("N100").Value = -3044.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1163.8182
$ws.Range("I113").Value = 1108.5
$ws.Range("J113").Value = 1311.3334
$ws.Range("K113").Value = 3325.5
$ws.Range("L113").Value = 3934.0002
$ws.Range("M113").Value = -1155.5
$ws.Range("N113").Value = -8274.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 49999.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 49999.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 49999.5
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -59919.5

Write-Output "Updated 33 leve rows across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."